$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F2:F4 (想去人数 column)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 996
$ws1.Range("F3").Value = 2043
$ws1.Range("F4").Value = 448

# Sheet "全部类型" - update F4:F6 (想去人数 column, same events duplicated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 996
$ws4.Range("F5").Value = 2043
$ws4.Range("F6").Value = 448
